# Fruta / hortaliza, semanal
# Insert a new weekly record at row 153 (pushing the existing 153..192
# block down by one row, to 154..193) and populate the new row with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 153; this shifts rows 153-192
# down to 154-193 (values, styles and formatting move with them).
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with this week's figures. The
# static descriptive columns match every other row in this block.
$ws.Range("A153").Value = 8
$ws.Range("B153").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C153").Value = 'Coquimbo'
$ws.Range("D153").Value = 44943
$ws.Range("E153").Value = 4
$ws.Range("F153").Value = 100112044
$ws.Range("G153").Value = 'Perejil'
$ws.Range("H153").Value = 'Sin especificar'
$ws.Range("I153").Value = 'Primera'
$ws.Range("J153").Value = 2000
$ws.Range("K153").Value = 3000
$ws.Range("L153").Value = 3500
$ws.Range("M153").Value = 3250
$ws.Range("N153").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O153").Value = 'Provincia del Elquí'
$ws.Range("P153").Value = 2167
$ws.Range("Q153").Value = 1.5
$ws.Range("R153").Value = 'Hortaliza'
